$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "edit1"
$ws.Range("B25").Value = "riya-morankar"
$ws.Range("C25").Value = "Merged"
$ws.Range("D25").Value = "compared ad merged"

# Force the date column to be stored as literal text ("2025-06-23"),
# matching the rest of the column, instead of being auto-parsed into a
# date serial number. Reset the number format back afterwards so the
# cell doesn't end up with a stray style compared to its neighbours.
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2025-06-23"
$ws.Range("E25").ClearFormats()

$ws.Range("F25").Value = "fee73b8b8fd20763dbba2cddf9bcd5df07ff197b"
